# Daily attendance processing - 2025-12-10 07:05:54
# Reorders the "Recorded By" (column G) comma-separated author lists so
# that any "System"/"system" entries come first, followed by the other
# entries, preserving their relative order. Cells without a System-like
# entry are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $value = $cell.Value2

    if ($value -eq $null) { continue }

    $parts = $value -split ',\s*'

    $systemParts = @()
    $otherParts = @()
    foreach ($part in $parts) {
        if ($part.ToLower() -eq 'system') {
            $systemParts += $part
        } else {
            $otherParts += $part
        }
    }

    if ($systemParts.Count -gt 0) {
        $newValue = ($systemParts + $otherParts) -join ', '
        if ($newValue -ne $value) {
            $cell.Value = $newValue
        }
    }
}
